$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 61.526261
$ws.Range("H2").Value = 184.578783
$ws.Range("I2").Value = 0.9684015388399598
$ws.Range("J2").Value = 0.9684015388399598
$ws.Range("M2").Value = 1.646992
$ws.Range("N2").Value = 4.940976
$ws.Range("O2").Value = 0.2071783517404009
$ws.Range("P2").Value = 0.2071783517404009
$ws.Range("Q2").Value = 101.333259656912
$ws.Range("R2").Value = 911.9993369122079
$ws.Range("S2").Value = 0.2006318346397307
$ws.Range("T2").Value = 0.2006318346397307
$ws.Range("G3").Value = 61.526261
$ws.Range("H3").Value = 184.578783
$ws.Range("I3").Value = 0.9684015388399598
$ws.Range("J3").Value = 0.9684015388399598
$ws.Range("O3").Value = 0.4685125322965616
$ws.Range("P3").Value = 0.4685125322965616
$ws.Range("Q3").Value = 229.154743673283
$ws.Range("R3").Value = 2062.392693059547
$ws.Range("S3").Value = 0.4537082572417966
$ws.Range("T3").Value = 0.4537082572417966
$ws.Range("G4").Value = 61.526261
$ws.Range("H4").Value = 184.578783
$ws.Range("I4").Value = 0.9684015388399598
$ws.Range("J4").Value = 0.9684015388399598
$ws.Range("M4").Value = 2.284352333333333
$ws.Range("N4").Value = 6.853057
$ws.Range("O4").Value = 0.2873531572796583
$ws.Range("P4").Value = 0.2873531572796583
$ws.Range("Q4").Value = 140.5476578766256
$ws.Range("R4").Value = 1264.928920889631
$ws.Range("S4").Value = 0.2782732397001421
$ws.Range("T4").Value = 0.2782732397001421
$ws.Range("G5").Value = 61.526261
$ws.Range("H5").Value = 184.578783
$ws.Range("I5").Value = 0.9684015388399598
$ws.Range("J5").Value = 0.9684015388399598
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2937863333333333
$ws.Range("N5").Value = 0.881359
$ws.Range("O5").Value = 0.03695595868337916
$ws.Range("P5").Value = 0.03695595868337916
$ws.Range("Q5").Value = 18.07557462289967
$ws.Range("R5").Value = 162.680171606097
$ws.Range("S5").Value = 0.03578820725829035
$ws.Range("T5").Value = 0.03578820725829036
$ws.Range("I6").Value = 0.002979850677668077
$ws.Range("J6").Value = 0.002979850677668078
$ws.Range("M6").Value = 1.646992
$ws.Range("N6").Value = 4.940976
$ws.Range("O6").Value = 0.2071783517404009
$ws.Range("P6").Value = 0.2071783517404009
$ws.Range("Q6").Value = 0.3118107214293334
$ws.Range("R6").Value = 2.806296492864
$ws.Range("S6").Value = 0.000617360551831789
$ws.Range("T6").Value = 0.0006173605518317891
$ws.Range("I7").Value = 0.002979850677668077
$ws.Range("J7").Value = 0.002979850677668078
$ws.Range("O7").Value = 0.4685125322965616
$ws.Range("P7").Value = 0.4685125322965616
$ws.Range("S7").Value = 0.001396097386859896
$ws.Range("T7").Value = 0.001396097386859896
$ws.Range("I8").Value = 0.002979850677668077
$ws.Range("J8").Value = 0.002979850677668078
$ws.Range("M8").Value = 2.284352333333333
$ws.Range("N8").Value = 6.853057
$ws.Range("O8").Value = 0.2873531572796583
$ws.Range("P8").Value = 0.2873531572796583
$ws.Range("Q8").Value = 0.4324766295497777
$ws.Range("R8").Value = 3.892289665948
$ws.Range("S8").Value = 0.0008562695004498513
$ws.Range("T8").Value = 0.0008562695004498514
$ws.Range("I9").Value = 0.002979850677668077
$ws.Range("J9").Value = 0.002979850677668078
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2937863333333333
$ws.Range("N9").Value = 0.881359
$ws.Range("O9").Value = 0.03695595868337916
$ws.Range("P9").Value = 0.03695595868337916
$ws.Range("Q9").Value = 0.05562002034177778
$ws.Range("R9").Value = 0.5005801830760001
$ws.Range("S9").Value = 0.0001101232385265408
$ws.Range("T9").Value = 0.0001101232385265409
$ws.Range("G10").Value = 0.6472316666666668
$ws.Range("H10").Value = 1.941695
$ws.Range("I10").Value = 0.01018719700821657
$ws.Range("J10").Value = 0.01018719700821657
$ws.Range("M10").Value = 1.646992
$ws.Range("N10").Value = 4.940976
$ws.Range("O10").Value = 0.2071783517404009
$ws.Range("P10").Value = 0.2071783517404009
$ws.Range("Q10").Value = 1.065985377146667
$ws.Range("R10").Value = 9.593868394320001
$ws.Range("S10").Value = 0.002110566685017053
$ws.Range("T10").Value = 0.002110566685017053
$ws.Range("G11").Value = 0.6472316666666668
$ws.Range("H11").Value = 1.941695
$ws.Range("I11").Value = 0.01018719700821657
$ws.Range("J11").Value = 0.01018719700821657
$ws.Range("O11").Value = 0.4685125322965616
$ws.Range("P11").Value = 0.4685125322965616
$ws.Range("Q11").Value = 2.410616284195001
$ws.Range("R11").Value = 21.69554655775501
$ws.Range("S11").Value = 0.004772829467323503
$ws.Range("T11").Value = 0.004772829467323503
$ws.Range("G12").Value = 0.6472316666666668
$ws.Range("H12").Value = 1.941695
$ws.Range("I12").Value = 0.01018719700821657
$ws.Range("J12").Value = 0.01018719700821657
$ws.Range("M12").Value = 2.284352333333333
$ws.Range("N12").Value = 6.853057
$ws.Range("O12").Value = 0.2873531572796583
$ws.Range("P12").Value = 0.2873531572796583
$ws.Range("Q12").Value = 1.478505167957222
$ws.Range("R12").Value = 13.306546511615
$ws.Range("S12").Value = 0.002927323224140921
$ws.Range("T12").Value = 0.002927323224140921
$ws.Range("G13").Value = 0.6472316666666668
$ws.Range("H13").Value = 1.941695
$ws.Range("I13").Value = 0.01018719700821657
$ws.Range("J13").Value = 0.01018719700821657
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2937863333333333
$ws.Range("N13").Value = 0.881359
$ws.Range("O13").Value = 0.03695595868337916
$ws.Range("P13").Value = 0.03695595868337916
$ws.Range("Q13").Value = 0.1901478181672222
$ws.Range("R13").Value = 1.711330363505
$ws.Range("S13").Value = 0.0003764776317350954
$ws.Range("T13").Value = 0.0003764776317350955
$ws.Range("G14").Value = 1.171018333333333
$ws.Range("H14").Value = 3.513055
$ws.Range("I14").Value = 0.01843141347415545
$ws.Range("J14").Value = 0.01843141347415545
$ws.Range("M14").Value = 1.646992
$ws.Range("N14").Value = 4.940976
$ws.Range("O14").Value = 0.2071783517404009
$ws.Range("P14").Value = 0.2071783517404009
$ws.Range("Q14").Value = 1.928657826853333
$ws.Range("R14").Value = 17.35792044168
$ws.Range("S14").Value = 0.003818589863821344
$ws.Range("T14").Value = 0.003818589863821344
$ws.Range("G15").Value = 1.171018333333333
$ws.Range("H15").Value = 3.513055
$ws.Range("I15").Value = 0.01843141347415545
$ws.Range("J15").Value = 0.01843141347415545
$ws.Range("O15").Value = 0.4685125322965616
$ws.Range("P15").Value = 0.4685125322965616
$ws.Range("Q15").Value = 4.361461295555
$ws.Range("R15").Value = 39.25315165999501
$ws.Range("S15").Value = 0.008635348200581537
$ws.Range("T15").Value = 0.008635348200581537
$ws.Range("G16").Value = 1.171018333333333
$ws.Range("H16").Value = 3.513055
$ws.Range("I16").Value = 0.01843141347415545
$ws.Range("J16").Value = 0.01843141347415545
$ws.Range("M16").Value = 2.284352333333333
$ws.Range("N16").Value = 6.853057
$ws.Range("O16").Value = 0.2873531572796583
$ws.Range("P16").Value = 0.2873531572796583
$ws.Range("Q16").Value = 2.675018462126111
$ws.Range("R16").Value = 24.075166159135
$ws.Range("S16").Value = 0.005296324854925404
$ws.Range("T16").Value = 0.005296324854925404
$ws.Range("G17").Value = 1.171018333333333
$ws.Range("H17").Value = 3.513055
$ws.Range("I17").Value = 0.01843141347415545
$ws.Range("J17").Value = 0.01843141347415545
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2937863333333333
$ws.Range("N17").Value = 0.881359
$ws.Range("O17").Value = 0.03695595868337916
$ws.Range("P17").Value = 0.03695595868337916
$ws.Range("Q17").Value = 0.3440291824161111
$ws.Range("R17").Value = 3.096262641745
$ws.Range("S17").Value = 0.0006811505548271668
$ws.Range("T17").Value = 0.0006811505548271669
